# Update excel.xlsx, run yarn start
# Inserts a new row into the "Snippets" table for the
# Excel.Workbook.getActiveShapeOrNullObject API snippet-metadata entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Snippets")
$tbl = $ws.ListObjects.Item(1)

# Insert a fresh worksheet row at row 336 (the row right after
# "getActiveCell" / before "getSelectedRanges"), shifting all the
# following data rows down by one.
$ws.Rows.Item(336).Insert()

# The table's defined range needs to grow by one row to cover the
# newly inserted row (A1:F389 -> A1:F390).
$tbl.Resize($ws.Range("A1:F390"))

# Populate the newly inserted row with the new snippet-metadata entry.
$ws.Range("A336").Value = "Excel"
$ws.Range("B336").Value = "Workbook"
$ws.Range("C336").Value = "getActiveShapeOrNullObject"
$ws.Range("D336").Value = 1
$ws.Range("E336").Value = "excel-shape-get-active"
$ws.Range("F336").Value = "getActiveShape"

# Move the active selection to where it ended up after the edit.
$ws.Range("F337").Select() | Out-Null
